$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$words = @(
    "leihen",
    "spielen",
    "schwören",
    "biegen",
    "stürzen",
    "wetten",
    "rasen",
    "trennen",
    "räumen",
    "äußern",
    "folgen",
    "fehlen",
    "regeln",
    "seufzen",
    "schlucken",
    "sammeln",
    "bleiben",
    "altern",
    "weichen",
    "lockern",
    "reizen",
    "backen",
    "warnen",
    "ächzen",
    "grüßen",
    "herzen",
    "sterben",
    "schmecken",
    "spüren",
    "wehtun",
    "zeugen",
    "stehlen",
    "zielen",
    "rufen",
    "schaden",
    "dürfen",
    "feiern",
    "trauen",
    "holen",
    "fügen",
    "bergen",
    "boxen",
    "wellen",
    "bluten",
    "saufen",
    "werden",
    "machen",
    "erben",
    "reisen",
    "strahlen",
    "geben",
    "schnellen",
    "prügeln",
    "morden",
    "achten",
    "ahnen",
    "ärgern",
    "schreiten",
    "jagen",
    "bellen",
    "danken",
    "ehren",
    "lügen",
    "mögen",
    "sprechen",
    "scheinen",
    "kämpfen",
    "heilen",
    "töten",
    "kriegen",
    "rechnen",
    "tropfen",
    "dringen",
    "schwingen",
    "rauchen",
    "schlagen",
    "landen",
    "kümmern",
    "scheitern",
    "helfen",
    "wenden",
    "arten",
    "leiden",
    "pissen",
    "feuern",
    "sprengen",
    "opfern",
    "segnen",
    "foltern",
    "drücken",
    "wundern",
    "drehen",
    "klingen",
    "fangen",
    "küssen",
    "fallen",
    "bauen",
    "hauen",
    "führen",
    "zeigen",
    "sinken",
    "tauchen",
    "gelten",
    "wüten",
    "stillen",
    "mauern",
    "malen",
    "dienen",
    "formen",
    "füllen",
    "schwächen",
    "betteln",
    "tollen",
    "klettern",
    "lieben",
    "kosten",
    "jubeln",
    "stimmen",
    "zünden",
    "irren",
    "weigern",
    "fließen",
    "steuern",
    "zahlen",
    "kichern",
    "schenken",
    "schätzen",
    "treiben",
    "stecken",
    "liefern",
    "runden",
    "knarren",
    "gründen",
    "lassen",
    "hören",
    "narren",
    "pfeifen",
    "schrecken",
    "ändern",
    "sichern",
    "heben",
    "sitzen",
    "heulen",
    "tanzen",
    "pflanzen",
    "loben",
    "dauern",
    "graben",
    "kehren",
    "enden",
    "brauchen",
    "schulden",
    "bitten",
    "bilden",
    "bieten",
    "liegen",
    "freuen",
    "fischen",
    "filmen",
    "wehren",
    "sorgen",
    "siegen",
    "fesseln",
    "münzen",
    "spinnen",
    "fällen",
    "handeln",
    "fühlen",
    "platzen",
    "flehen",
    "planen",
    "sperren",
    "decken",
    "wirken",
    "fragen",
    "quälen",
    "flüchten",
    "binden",
    "wachsen",
    "suchen",
    "lesen",
    "fahren",
    "streichen",
    "leeren",
    "stürmen",
    "greifen",
    "werfen",
    "meistern",
    "zögern",
    "nerven",
    "stammen",
    "melden"
)

for ($i = 0; $i -lt $words.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $words[$i]
}

